# Apply the "electric rice cooker" sheet edit:
#  1. Swap the "B" and "C" quarter rows (columns A:E) within each year block.
#  2. Remove columns F (产销率) and G (销售量) entirely, including the header cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A scratch row well below the used data (data only goes to row 65) used as a
# temporary holding area so Copy (which preserves blank-cell identity, unlike
# plain Value2 array assignment) can perform a 3-way row swap. Each
# destination range is explicitly cleared before being pasted into, because
# Copy-ing a genuinely blank source cell onto a populated destination cell
# leaves the destination's old value untouched instead of blanking it.
$scratchRow = 200

# --- Step 1: swap B/C quarter-row data (columns A:E) for every year block ---
# Data starts at row 2 (2000年A) in 4-row blocks: A, B, C, D.
# The "B" row of block n is at row (3 + 4*n); the "C" row is at row (4 + 4*n).
for ($yearIndex = 0; $yearIndex -lt 16; $yearIndex++) {
    $bRow = 3 + (4 * $yearIndex)
    $cRow = 4 + (4 * $yearIndex)

    $rangeB = $ws.Range("A$bRow`:E$bRow")
    $rangeC = $ws.Range("A$cRow`:E$cRow")
    $rangeScratch = $ws.Range("A$scratchRow`:E$scratchRow")

    $rangeB.Copy($rangeScratch)
    $rangeB.ClearContents()
    $rangeC.Copy($rangeB)
    $rangeC.ClearContents()
    $rangeScratch.Copy($rangeC)
    $rangeScratch.Clear()
}

# --- Step 2: delete columns F and G (header + data) ---
$ws.Range("F1:G65").Delete()

Write-Host "Edit applied."
